$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V2 Command")

# Insert a new row before row 31 (shifts old row31.. down to row32..)
$ws.Rows.Item(31).Insert()

$ws.Cells.Item(31, 2).Value = 16            # B31
$ws.Cells.Item(31, 3).Value = "HaiLzd Servo Command"   # C31
$ws.Cells.Item(31, 4).Value = "V2"          # D31 (same as v70 string, placeholder - fixed below)
$ws.Cells.Item(31, 5).Value = "{cmd} {parm….}"          # E31
$ws.Cells.Item(31, 6).Value = "A9 9A 03 16 01 1A ED"    # F31
